$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1049.875
$ws.Range("I2").Value = 985.5714
$ws.Range("K2").Value = 985.5714
$ws.Range("M2").Value = -872.5714
$ws.Range("H43").Value = 7444.591
$ws.Range("I43").Value = 6080
$ws.Range("K43").Value = 6080
$ws.Range("M43").Value = -6011
$ws.Range("H98").Value = 2799.9412
$ws.Range("I98").Value = 3292.8572
$ws.Range("K98").Value = 3292.8572
$ws.Range("M98").Value = -1794.8572
$ws.Range("H122").Value = 2799.9412
$ws.Range("I122").Value = 3292.8572
$ws.Range("K122").Value = 9878.571599999999
$ws.Range("M122").Value = -7428.571599999999
$ws.Range("H131").Value = 2095.5
$ws.Range("I131").Value = 2095.5
$ws.Range("K131").Value = 6286.5
$ws.Range("M131").Value = -1246.5
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H138").Value = 3429.34
$ws.Range("I138").Value = 2236.3333
$ws.Range("J138").Value = 3746.4683
$ws.Range("K138").Value = 6708.999899999999
$ws.Range("L138").Value = 11239.4049
$ws.Range("M138").Value = -1568.999899999999
$ws.Range("N138").Value = -21519.4049
$ws.Range("H141").Value = 7718.5
$ws.Range("I141").Value = 7718.5
$ws.Range("K141").Value = 23155.5
$ws.Range("M141").Value = -17975.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15133998
$ws.Range("I32").Value = 16381925
$ws.Range("K32").Value = 16381925
$ws.Range("M32").Value = -16381638
$ws.Range("H45").Value = 5150.25
$ws.Range("I45").Value = 4880.4
$ws.Range("K45").Value = 4880.4
$ws.Range("M45").Value = -4503.4
$ws.Range("H102").Value = 2352.0833
$ws.Range("I102").Value = 1802.8889
$ws.Range("K102").Value = 1802.8889
$ws.Range("M102").Value = -180.8888999999999
$ws.Range("H122").Value = 4443.9287
$ws.Range("I122").Value = 3326.1667
$ws.Range("J122").Value = 5282.25
$ws.Range("K122").Value = 9978.500100000001
$ws.Range("L122").Value = 15846.75
$ws.Range("M122").Value = -7528.500100000001
$ws.Range("N122").Value = -20746.75
$ws.Range("H137").Value = 200000
$ws.Range("J137").Value = 200000
$ws.Range("L137").Value = 200000
$ws.Range("N137").Value = -210200

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 3757
$ws.Range("I16").Value = 3700
$ws.Range("J16").Value = 3785.5
$ws.Range("K16").Value = 3700
$ws.Range("L16").Value = 3785.5
$ws.Range("M16").Value = -3530
$ws.Range("N16").Value = -4125.5
$ws.Range("H40").Value = 58378
$ws.Range("J40").Value = 58378
$ws.Range("L40").Value = 58378
$ws.Range("N40").Value = -58908
$ws.Range("H86").Value = 1520.9286
$ws.Range("I86").Value = 1237.9
$ws.Range("J86").Value = 2228.5
$ws.Range("K86").Value = 1237.9
$ws.Range("L86").Value = 2228.5
$ws.Range("M86").Value = -114.9000000000001
$ws.Range("N86").Value = -4474.5
$ws.Range("H87").Value = 164000
$ws.Range("J87").Value = 164000
$ws.Range("L87").Value = 164000
$ws.Range("N87").Value = -166496
$ws.Range("H89").Value = 1520.9286
$ws.Range("I89").Value = 1237.9
$ws.Range("J89").Value = 2228.5
$ws.Range("K89").Value = 6189.5
$ws.Range("L89").Value = 11142.5
$ws.Range("M89").Value = -573.5
$ws.Range("N89").Value = -22374.5
$ws.Range("H90").Value = 164000
$ws.Range("J90").Value = 164000
$ws.Range("L90").Value = 492000
$ws.Range("N90").Value = -504480
$ws.Range("H96").Value = 36115
$ws.Range("I96").Value = 4175
$ws.Range("J96").Value = 99995
$ws.Range("K96").Value = 4175
$ws.Range("L96").Value = 99995
$ws.Range("M96").Value = -1429
$ws.Range("N96").Value = -105487
$ws.Range("H105").Value = 2500.2778
$ws.Range("I105").Value = 2187.8125
$ws.Range("K105").Value = 2187.8125
$ws.Range("M105").Value = -440.8125

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1854.1091
$ws.Range("I31").Value = 1342.3478
$ws.Range("J31").Value = 2221.9375
$ws.Range("K31").Value = 1342.3478
$ws.Range("L31").Value = 2221.9375
$ws.Range("M31").Value = -1047.3478
$ws.Range("N31").Value = -2811.9375
$ws.Range("H34").Value = 1854.1091
$ws.Range("I34").Value = 1342.3478
$ws.Range("J34").Value = 2221.9375
$ws.Range("K34").Value = 1342.3478
$ws.Range("L34").Value = 2221.9375
$ws.Range("M34").Value = -1140.3478
$ws.Range("N34").Value = -2625.9375
$ws.Range("H58").Value = 4470.2856
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H99").Value = 3534.1428
$ws.Range("I99").Value = 3213
$ws.Range("K99").Value = 3213
$ws.Range("M99").Value = -1715
$ws.Range("H105").Value = 2208.7273
$ws.Range("I105").Value = 2099.625
$ws.Range("J105").Value = 2499.6667
$ws.Range("K105").Value = 2099.625
$ws.Range("L105").Value = 2499.6667
$ws.Range("M105").Value = -352.625
$ws.Range("N105").Value = -5993.6667
$ws.Range("H126").Value = 3534.1428
$ws.Range("I126").Value = 3213
$ws.Range("K126").Value = 9639
$ws.Range("M126").Value = -7169
$ws.Range("H132").Value = 4554
$ws.Range("I132").Value = 3698
$ws.Range("K132").Value = 11094
$ws.Range("M132").Value = -8564
$ws.Range("H134").Value = 2000
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465
$ws.Range("H136").Value = 4470.2856
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 606.8
$ws.Range("I44").Value = 13
$ws.Range("K44").Value = 39
$ws.Range("M44").Value = 359
$ws.Range("H68").Value = 2099.484
$ws.Range("J68").Value = 2269.1738
$ws.Range("L68").Value = 6807.5214
$ws.Range("N68").Value = -8429.5214
$ws.Range("H71").Value = 2099.484
$ws.Range("J71").Value = 2269.1738
$ws.Range("L71").Value = 20422.5642
$ws.Range("N71").Value = -28534.5642
$ws.Range("H107").Value = 1214.7273
$ws.Range("J107").Value = 1562.8889
$ws.Range("L107").Value = 4688.6667
$ws.Range("N107").Value = -8528.6667
$ws.Range("H132").Value = 2267.88
$ws.Range("J132").Value = 2194.5789
$ws.Range("L132").Value = 19751.2101
$ws.Range("N132").Value = -24811.2101

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2387.6428
$ws.Range("I22").Value = 1832.4445
$ws.Range("K22").Value = 1832.4445
$ws.Range("M22").Value = -1537.4445
$ws.Range("H27").Value = 2387.6428
$ws.Range("I27").Value = 1832.4445
$ws.Range("K27").Value = 1832.4445
$ws.Range("M27").Value = -1725.4445
$ws.Range("H132").Value = 9693.125
$ws.Range("I132").Value = 10577.143
$ws.Range("K132").Value = 31731.429
$ws.Range("M132").Value = -29201.429
$ws.Range("H137").Value = 117499.5
$ws.Range("J137").Value = 117499.5
$ws.Range("L137").Value = 117499.5
$ws.Range("N137").Value = -127699.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5302.3335
$ws.Range("I62").Value = 5092.778
$ws.Range("J62").Value = 5616.6665
$ws.Range("K62").Value = 5092.778
$ws.Range("L62").Value = 5616.6665
$ws.Range("M62").Value = -4468.778
$ws.Range("N62").Value = -6864.6665
$ws.Range("H65").Value = 5302.3335
$ws.Range("I65").Value = 5092.778
$ws.Range("J65").Value = 5616.6665
$ws.Range("K65").Value = 25463.89
$ws.Range("L65").Value = 28083.3325
$ws.Range("M65").Value = -22343.89
$ws.Range("N65").Value = -34323.3325
$ws.Range("H122").Value = 58828750
$ws.Range("I122").Value = 62505252
$ws.Range("K122").Value = 187515756
$ws.Range("M122").Value = -187513306
$ws.Range("H126").Value = 4346
$ws.Range("I126").Value = 4269.1177
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 12807.3531
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -10337.3531
$ws.Range("N126").Value = -19938.5
$ws.Range("H131").Value = 128997.86
$ws.Range("J131").Value = 128997.86
$ws.Range("L131").Value = 128997.86
$ws.Range("N131").Value = -139077.86
$ws.Range("H132").Value = 1865.1111
$ws.Range("I132").Value = 1760.2858
$ws.Range("K132").Value = 5280.857400000001
$ws.Range("M132").Value = -2750.857400000001
$ws.Range("H136").Value = 31342.686
$ws.Range("I136").Value = 2327.4
$ws.Range("K136").Value = 6982.200000000001
$ws.Range("M136").Value = -4432.200000000001
